# adapting workflow to percentage load changing and load increasing by excel
$wb = $excel.ActiveWorkbook

# --- "Coupling Parameters" sheet ---
$ws = $wb.Worksheets.Item("Coupling Parameters")

# B3: look-ahead formula 2031 -> 2025
$ws.Range("B3").Formula = "=2025"

# B24: dummy_capacity_to_be_installed 2000 -> 2500
$ws.Range("B24").Value = 2500

# B31: increasingLoad_representativeYear_Excel, "None" -> name of excel used for load increase
$ws.Range("B31").Value = "2020-2050_basedon2004.xlsx"

# C31: updated helper text for B31
$ws.Range("C31").Value = "None or the name of the excel that future load is based on "

# --- "optionsConfig" sheet ---
$ws3 = $wb.Worksheets.Item("optionsConfig")

# G6: new allowed option for the dropdown list backing B31/B32
$ws3.Range("G6").Value = "2020-2050_basedon2004.xlsx"

# --- Selections (match final cursor position per sheet) ---
$ws3.Activate()
$ws3.Range("I14").Select()

$ws.Activate()
$ws.Range("C10").Select()
